$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '65.620.73'
$c.ClearFormats()
$ws.Range('E2').Value = '  -0.82%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.670.09'
$c.ClearFormats()
$ws.Range('E3').Value = '  -0.83%  '
$ws.Range('E4').Value = '  +0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '600.92'
$c.ClearFormats()
$ws.Range('E5').Value = '  -1.54%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '156.96'
$c.ClearFormats()
$ws.Range('E6').Value = '  -1.49%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +4.62%  '
$ws.Range('E9').Value = '  +3.05%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.403'
$c.ClearFormats()
$ws.Range('E10').Value = '  -0.79%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '5.87'
$c.ClearFormats()
$ws.Range('E11').Value = '  -2.45%  '
$ws.Range('E12').Value = '  -0.45%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '29.33'
$c.ClearFormats()
$ws.Range('E13').Value = '  -3.37%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '0.0000200'
$c.ClearFormats()
$ws.Range('E14').Value = '  -6.38%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.149.76'
$c.ClearFormats()
$ws.Range('E15').Value = '  -0.84%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '65.466.32'
$c.ClearFormats()
$ws.Range('E16').Value = '  -0.83%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '2.669.20'
$c.ClearFormats()
$ws.Range('E17').Value = '  -0.88%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '12.81'
$c.ClearFormats()
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('E19').Value = '  -2.35%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '7.65'
$c.ClearFormats()
$ws.Range('E20').Value = '  +1.53%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '350.18'
$c.ClearFormats()
$ws.Range('E21').Value = '  -3.61%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E22').Value = '  -0.07%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '69.75'
$c.ClearFormats()
$ws.Range('E23').Value = '  -0.60%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '0.0000111'
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.21%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '9.72'
$c.ClearFormats()
$ws.Range('E25').Value = '  -0.62%  '
$ws.Range('E26').Value = '  -2.13%  '
$ws.Range('E27').Value = '  -3.74%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.59'
$c.ClearFormats()
$ws.Range('E28').Value = '  -6.66%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '8.07'
$c.ClearFormats()
$ws.Range('E29').Value = '  -2.78%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.ClearFormats()
$ws.Range('E30').Value = '  +0.28%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '2.15'
$c.ClearFormats()
$ws.Range('E31').Value = '  -3.07%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '532.45'
$c.ClearFormats()
$ws.Range('E32').Value = '  -0.68%  '
$ws.Range('E33').Value = '  -2.61%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '6.49'
$c.ClearFormats()
$ws.Range('E34').Value = '  -2.13%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '5.51'
$c.ClearFormats()
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('E36').Value = '  -2.42%  '
$ws.Range('E37').Value = '  -2.03%  '
$ws.Range('E38').Value = '  +0.01%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '158.48'
$c.ClearFormats()
$ws.Range('E39').Value = '  -2.81%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '1.93'
$c.ClearFormats()
$ws.Range('E40').Value = '  -3.66%  '
$ws.Range('E41').Value = '  +0.01%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '42.73'
$c.ClearFormats()
$ws.Range('E42').Value = '  -0.55%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '164.60'
$c.ClearFormats()
$ws.Range('E43').Value = '  -3.79%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '4.14'
$c.ClearFormats()
$ws.Range('E44').Value = '  -1.38%  '
$ws.Range('E45').Value = '  -0.66%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0612'
$c.ClearFormats()
$ws.Range('E46').Value = '  -0.88%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '22.98'
$c.ClearFormats()
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('E48').Value = '  -2.70%  '
$ws.Range('E49').Value = '  -2.73%  '
$ws.Range('E50').Value = '  +1.91%  '
$ws.Range('E51').Value = '  -1.16%  '
